$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.571.14"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.628.82"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "593.93"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "168.13"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "2.629.40"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "27.66"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "3.117.61"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "67.791.77"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "2.626.15"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "12.03"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").Value = "8.05"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "358.14"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "10.36"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("D27").Value = "70.02"
$ws.Range("D28").Value = "2.766.47"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "548.01"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("D38").Value = "158.16"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").Value = "19.07"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "18.29"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "152.88"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "0.581"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "3.81"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "0.0772"
$ws.Range("E51").Value = "  -0.80%  "
